$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Dillon Brooks", "SG,SF", "Houston Rockets"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
